# Auto-generated edit script: update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.262.26'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '3.507.71'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.43%  '
$ws.Range("D7").Value = '3.508.57'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.375'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '4.107.69'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000180'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '3.518.04'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '64.305.80'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.569'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = '3.650.03'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  +4.35%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +1.07%  '
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("D34").Value = '3.525.99'
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.55'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '164.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("E42").Value = '  -2.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.809'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.82%  '
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").Value = '2.474.72'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.919'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.55%  '
